$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.233.49"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.519.52"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.50"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.80"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "2.540.28"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.60"
$ws.Range("E12").Value = "  +4.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "2.964.59"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.65"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "59.153.46"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "2.532.79"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.56"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.02"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  +3.60%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  +6.64%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.52"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.72"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.93"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.70"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.35"
$ws.Range("E42").Value = "  -6.12%  "
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.76"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.27"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("E51").Value = "  -2.33%  "
